# Append a newly-scraped Lancers listing to the "ランサーズ" sheet.
# A fresh row is inserted right after the two untouched rows (2-3); the
# rows that used to follow shift down by one and get their "fetched at"
# timestamp refreshed to the new scrape time, matching a re-run of the
# scraper that prepends its newest find ahead of previously-seen ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-12-29 12:40:10"

# Widen column D (price) slightly to fit the larger new price range text.
# (ColumnWidth round-trips through this host's char-width conversion with a
# constant +5/6 offset vs. the stored OOXML width, so back it out here to
# land on exactly width="32".)
$ws.Columns.Item(4).ColumnWidth = 32 - 5/6

# Shift rows 4 (and below) down by one to make room for the new listing.
$ws.Rows.Item(4).Insert()

# Refresh the scrape timestamp on every row (new row included below).
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp

# Populate the newly inserted row 4 with the new listing's data.
$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "【SNSアプリ開発】AndroidとiOSのインスタグラム風アプリ制作依頼"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5462964"
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = "◆開発 ◇アプリ"

# The row-insert doesn't renumber the worksheet's hyperlink list, so rebuild
# it from scratch against the now-correct row positions.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5427956") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5454210") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5462964") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5462891") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5462712") | Out-Null

# Hyperlinks.Add stamps its own style variant; pin the URL column back to
# the workbook's shared "Hyperlink" cell style used by the other rows.
$ws.Range("F2:F6").Style = "Hyperlink"
